# Apply edits described by the diff:
#  - C3 changes from duplicate "K086-LM" to "M254-VM"
#  - D3 changes from duplicate "K086-LM" to "M199-DS"
#  - C9 changes from duplicate "K137-DT" to "M198-CN"
#  - selection/active cell moves to E5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = "M254-VM"
$ws.Range("D3").Value = "M199-DS"
$ws.Range("C9").Value = "M198-CN"

$ws.Range("E5").Select()
